# Update workbook to reflect the published CDA FHIR logical model
# refresh (patches #241): new matchbox-patch version, new generation
# date, corrected Contact, corrected Binding Value Set URL for
# Entity.code, and corrected Min / Base Min cardinality for
# Entity.classCode (0..1 instead of 1..1).

$wb = $excel.ActiveWorkbook

# --- "Metadata" sheet -----------------------------------------------
$wsMeta = $wb.Worksheets.Item("Metadata")

# Version
$wsMeta.Range("B3").Value = "2.0.0-sd-202406-matchbox-patch"

# Date
$wsMeta.Range("B8").Value = "2024-06-19T17:47:42+02:00"

# Contact
$wsMeta.Range("B10").Value = "HL7 International - Structured Documents (http://www.hl7.org/Special/committees/structure, structdog@lists.HL7.org)"

# --- "Elements" sheet -------------------------------------------------
$wsElem = $wb.Worksheets.Item("Elements")

# Entity.code Binding Value Set
$wsElem.Range("Z15").Value = "http://hl7.org/cda/stds/core/ValueSet/CDAEntityCode"

# Entity.classCode Min (F12) and Base Min (AG12): 1 -> 0
# Force text formatting first so the value is stored as the text "0"
# (matching Min/Max text convention used throughout this table)
# rather than being auto-coerced into a numeric value.
$wsElem.Range("F12").NumberFormat = "@"
$wsElem.Range("F12").Value = "0"

$wsElem.Range("AG12").NumberFormat = "@"
$wsElem.Range("AG12").Value = "0"
